# menambahkan model 2 dan 3
# Splits the old "Development Front-End Web" task row into four new rows
# (image-generator + 3 model-testing rows), pushing the remaining Waiting
# tasks (Configurasi Flask / Development Back-End Web / Deployment) down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. make room: insert 3 blank rows right before the old row 11 ----------
$ws.Range("A11:A13").EntireRow.Insert()

# After this insert:
#   row 10            -> unchanged ("Save Model 1 dan 2")
#   rows 11,12,13      -> brand new, blank
#   old row 11 (was "Development Front-End Web") -> now row 14
#   old row 12 (was "Configurasi Flask")          -> now row 15
#   old row 13 (was "Development Back-End Web")   -> now row 16
#   old row 14 (was "Deployment ")                -> now row 17

# --- 2. carry formatting onto the new blank rows ----------------------------
# Row 11 continues the "Done" block styling, like row 10 above it.
$ws.Range("A10:F10").Copy()
$ws.Range("A11:F11").PasteSpecial(-4122)

# Rows 12 & 13 continue the styling that belonged to the task block that
# starts at (now) row 14.
$ws.Range("A14:F14").Copy()
$ws.Range("A12:F13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row heights for the newly inserted rows match the rest of the table.
$ws.Rows.Item(11).RowHeight = 15.5
$ws.Rows.Item(12).RowHeight = 15.5
$ws.Rows.Item(13).RowHeight = 15.5

# --- 3. fill in the new job rows (11-14) ------------------------------------
$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = "Membuat image generator"
$ws.Cells.Item(11, 3).Value = "8 Nopember"
$ws.Cells.Item(11, 4).Value = "8 Nopember"
$ws.Cells.Item(11, 5).Value = "Husein"
$ws.Cells.Item(11, 6).Value = "Done"

$ws.Cells.Item(12, 1).Value = 9
$ws.Cells.Item(12, 2).Value = "Testing Model 1"
$ws.Cells.Item(12, 3).Value = "8 Nopember"
$ws.Cells.Item(12, 4).Value = "8 Nopember"
$ws.Cells.Item(12, 5).Value = "Husein"
$ws.Cells.Item(12, 6).Value = "Done"

$ws.Cells.Item(13, 1).Value = 10
$ws.Cells.Item(13, 2).Value = "Testing Model 2"
$ws.Cells.Item(13, 3).Value = "9 Nopember"
$ws.Cells.Item(13, 4).Value = "9 Nopember"
$ws.Cells.Item(13, 5).Value = "Adhi"
$ws.Cells.Item(13, 6).Value = "Done"

$ws.Cells.Item(14, 1).Value = 11
$ws.Cells.Item(14, 2).Value = "Testing Model 3"
$ws.Cells.Item(14, 3).Value = "9 Nopember"
$ws.Cells.Item(14, 4).Value = "9 Nopember"
$ws.Cells.Item(14, 5).Value = "Adhi"
$ws.Cells.Item(14, 6).Value = "Done"

# --- 4. renumber the remaining Waiting rows that shifted down ---------------
$ws.Cells.Item(15, 1).Value = 12          # Configurasi Flask
$ws.Cells.Item(16, 1).Value = 13          # Development Back-End Web
$ws.Cells.Item(17, 1).Value = 14          # Deployment

# --- 5. view/selection bookkeeping ------------------------------------------
$ws.Range("C16").Select()
